$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 216, pushing existing rows 216-250 down to 218-252.
$ws.Rows.Item(216).Resize(2, 1).Insert()

# Fill in the two newly inserted rows (216 and 217) with their data.
# Row 216
$ws.Cells.Item(216, 1).Value2 = 3
$ws.Cells.Item(216, 2).Value2 = "Femacal de La Calera"
$ws.Cells.Item(216, 3).Value2 = "Coquimbo"
$ws.Cells.Item(216, 4).Value2 = 44522
$ws.Cells.Item(216, 4).NumberFormat = $ws.Cells.Item(218, 4).NumberFormat
$ws.Cells.Item(216, 5).Value2 = 5
$ws.Cells.Item(216, 6).Value2 = 100112031
$ws.Cells.Item(216, 7).Value2 = "Poroto verde"
$ws.Cells.Item(216, 8).Value2 = "Magnum"
$ws.Cells.Item(216, 9).Value2 = "Primera"
$ws.Cells.Item(216, 10).Value2 = 70
$ws.Cells.Item(216, 11).Value2 = 40000
$ws.Cells.Item(216, 12).Value2 = 41000
$ws.Cells.Item(216, 13).Value2 = 40500
$ws.Cells.Item(216, 14).Value2 = "$/malla 25 kilos"
$ws.Cells.Item(216, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(216, 16).Value2 = 1620
$ws.Cells.Item(216, 17).Value2 = 25
$ws.Cells.Item(216, 18).Value2 = "Hortaliza"

# Row 217
$ws.Cells.Item(217, 1).Value2 = 3
$ws.Cells.Item(217, 2).Value2 = "Femacal de La Calera"
$ws.Cells.Item(217, 3).Value2 = "Coquimbo"
$ws.Cells.Item(217, 4).Value2 = 44522
$ws.Cells.Item(217, 4).NumberFormat = $ws.Cells.Item(218, 4).NumberFormat
$ws.Cells.Item(217, 5).Value2 = 5
$ws.Cells.Item(217, 6).Value2 = 100112031
$ws.Cells.Item(217, 7).Value2 = "Poroto verde"
$ws.Cells.Item(217, 8).Value2 = "Magnum"
$ws.Cells.Item(217, 9).Value2 = "Segunda"
$ws.Cells.Item(217, 10).Value2 = 35
$ws.Cells.Item(217, 11).Value2 = 34000
$ws.Cells.Item(217, 12).Value2 = 34000
$ws.Cells.Item(217, 13).Value2 = 34000
$ws.Cells.Item(217, 14).Value2 = "$/malla 25 kilos"
$ws.Cells.Item(217, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(217, 16).Value2 = 1360
$ws.Cells.Item(217, 17).Value2 = 25
$ws.Cells.Item(217, 18).Value2 = "Hortaliza"
